$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo: X7 "Mutliple" -> "Multiple"
$ws.Range("X7").Value = "Multiple"

# Fill all blank cells in A2:AA9 with the literal text "null"
$rng = $ws.Range("A2:AA9")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq $null) {
        $cell.NumberFormat = "@"
        $cell.Value = "null"
    }
}

# Add conditional formatting rule: highlight blank cells (after TRIM) in red across A2:AA9
$cfRange = $ws.Range("A2:AA9")
$cf = $cfRange.FormatConditions.Add(2, 0, "LEN(TRIM(A2))=0")
$cf.Interior.Color = 255
$cf.SetFirstPriority()

# Update sheet view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("Z25").Select()
